$wb = $excel.ActiveWorkbook

# --- Update PipeDataXlIn (input) sheet literal values ---
$wsIn = $wb.Worksheets.Item("PipeDataXlIn")

$wsIn.Range("E2").Value = 0.019999999999999997
$wsIn.Range("G2").Value = 116.19691664343739
$wsIn.Range("H2").Value = 112.50814631914275

$wsIn.Range("E3").Value = 0.0089999999999999143
$wsIn.Range("G3").Value = 111.91814631914275
$wsIn.Range("H3").Value = 111.23466110510971

$wsIn.Range("E4").Value = 0.0056000000000000216
$wsIn.Range("G4").Value = 110.64466110510971

# --- Update selection on PipeDataXlOut (output) sheet ---
$wsOut = $wb.Worksheets.Item("PipeDataXlOut")
$wsOut.Activate()
$wsOut.Range("L7").Select()
